# Generate Report for Handoff
#
# Updates the localization-status report to reflect a fresh handoff:
#  - Status cells move from "Handed back: in sync with en-US" to
#    "Ready for handoff"
#  - The "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
#    timestamps are refreshed
#  - The (now shorter) status text no longer needs the wide columns that
#    used to hold the long "Handed back..." message, so those columns are
#    narrowed back down

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff"
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# --- Refreshed timestamps
$overview.Range("G2").Value = "2016-08-12 07:13:25"
$zhcn.Range("H2").Value = "2016-08-12 07:13:19"
$dede.Range("H2").Value = "2016-08-12 07:13:25"

# --- Narrow the status columns now that the status text is shorter.
# ColumnWidth is expressed in (quantized) Excel character-width units;
# 16.38265482584637 is the input that lands on the target display width.
$overview.Columns.Item(5).ColumnWidth = 16.38265482584637
$overview.Columns.Item(6).ColumnWidth = 16.38265482584637
$zhcn.Columns.Item(3).ColumnWidth = 16.38265482584637
$dede.Columns.Item(3).ColumnWidth = 16.38265482584637
